$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 111112320
$ws.Range("I18").Value = 1360
$ws.Range("K18").Value = 1360
$ws.Range("M18").Value = -1076

# Row 49
$ws.Range("H49").Value = 14999.75
$ws.Range("I49").Value = 14999
$ws.Range("K49").Value = 44997
$ws.Range("M49").Value = -44861

# Row 52
$ws.Range("H52").Value = 437
$ws.Range("I52").Value = 465.66666
$ws.Range("J52").Value = 351
$ws.Range("K52").Value = 1396.99998
$ws.Range("L52").Value = 1053
$ws.Range("M52").Value = -1236.99998
$ws.Range("N52").Value = -1373

# Row 74
$ws.Range("H74").Value = 15051.826
$ws.Range("I74").Value = 15808.471
$ws.Range("J74").Value = 12908
$ws.Range("K74").Value = 15808.471
$ws.Range("L74").Value = 12908
$ws.Range("M74").Value = -14872.471
$ws.Range("N74").Value = -14780

# Row 77
$ws.Range("H77").Value = 15051.826
$ws.Range("I77").Value = 15808.471
$ws.Range("J77").Value = 12908
$ws.Range("K77").Value = 79042.355
$ws.Range("L77").Value = 64540
$ws.Range("M77").Value = -74362.355
$ws.Range("N77").Value = -73900

# Row 98
$ws.Range("H98").Value = 4099.846
$ws.Range("I98").Value = 4099.846
$ws.Range("K98").Value = 4099.846
$ws.Range("M98").Value = -2601.846

# Row 113
$ws.Range("H113").Value = 3501
$ws.Range("I113").Value = 3335.3333
$ws.Range("K113").Value = 3335.3333
$ws.Range("M113").Value = -81.33329999999978

# Row 122
$ws.Range("H122").Value = 4099.846
$ws.Range("I122").Value = 4099.846
$ws.Range("K122").Value = 12299.538
$ws.Range("M122").Value = -9849.537999999999

# Row 125
$ws.Range("H125").Value = 1215
$ws.Range("J125").Value = 1215
$ws.Range("L125").Value = 10935
$ws.Range("N125").Value = -15855

# Row 132
$ws.Range("H132").Value = 5983.7837
$ws.Range("I132").Value = 4761.1387
$ws.Range("J132").Value = 49999
$ws.Range("K132").Value = 14283.4161
$ws.Range("L132").Value = 149997
$ws.Range("M132").Value = -11753.4161
$ws.Range("N132").Value = -155057

$ws = $wb.Sheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3108.1636
$ws.Range("I32").Value = 2885.3962
$ws.Range("J32").Value = 9011.5
$ws.Range("K32").Value = 2885.3962
$ws.Range("L32").Value = 9011.5
$ws.Range("M32").Value = -2598.3962
$ws.Range("N32").Value = -9585.5

# Row 74
$ws.Range("H74").Value = 244424.61
$ws.Range("I74").Value = 506904.8
$ws.Range("K74").Value = 506904.8
$ws.Range("M74").Value = -506030.8

# Row 77
$ws.Range("H77").Value = 244424.61
$ws.Range("I77").Value = 506904.8
$ws.Range("K77").Value = 2534524
$ws.Range("M77").Value = -2530156

# Row 132
$ws.Range("H132").Value = 3180.2856
$ws.Range("I132").Value = 1852.4667
$ws.Range("K132").Value = 5557.4001
$ws.Range("M132").Value = -3027.4001

$ws = $wb.Sheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4128.364
$ws.Range("I134").Value = 3906.7637
$ws.Range("J134").Value = 5236.364
$ws.Range("K134").Value = 11720.2911
$ws.Range("L134").Value = 15709.092
$ws.Range("M134").Value = -9185.2911
$ws.Range("N134").Value = -20779.092

$ws = $wb.Sheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3978.3333
$ws.Range("I31").Value = 2862.6072
$ws.Range("K31").Value = 2862.6072
$ws.Range("M31").Value = -2567.6072

# Row 34
$ws.Range("H34").Value = 3978.3333
$ws.Range("I34").Value = 2862.6072
$ws.Range("K34").Value = 2862.6072
$ws.Range("M34").Value = -2660.6072

# Row 103
$ws.Range("H103").Value = 10538.857
$ws.Range("I103").Value = 10538.857
$ws.Range("K103").Value = 10538.857
$ws.Range("M103").Value = -9366.857

# Row 132
$ws.Range("H132").Value = 3491.24
$ws.Range("I132").Value = 3506.0908
$ws.Range("J132").Value = 3479.5715
$ws.Range("K132").Value = 10518.2724
$ws.Range("L132").Value = 10438.7145
$ws.Range("M132").Value = -7988.2724
$ws.Range("N132").Value = -15498.7145

# Row 139
$ws.Range("H139").Value = 98998.5
$ws.Range("J139").Value = 98998.5
$ws.Range("L139").Value = 98998.5
$ws.Range("N139").Value = -109278.5

$ws = $wb.Sheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 847.5
$ws.Range("I92").Value = 795
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 2385
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196
$ws.Range("M92").Value = -1137

# Row 134
$ws.Range("H134").Value = 2262.7144
$ws.Range("I134").Value = 2262.7144
$ws.Range("K134").Value = 6788.1432
$ws.Range("M134").Value = -1718.1432

$ws = $wb.Sheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1827.8572
$ws.Range("J102").Value = 2920
$ws.Range("L102").Value = 2920
$ws.Range("N102").Value = -6164

# Row 122
$ws.Range("H122").Value = 6520.5625
$ws.Range("I122").Value = 5236.9
$ws.Range("J122").Value = 8660
$ws.Range("K122").Value = 15710.7
$ws.Range("L122").Value = 25980
$ws.Range("M122").Value = -13260.7
$ws.Range("N122").Value = -30880

$ws = $wb.Sheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1954.9
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1954.9
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1954.9
$ws.Range("N46").Value = -2330.9
$ws.Range("M46").ClearContents()

# Row 56
$ws.Range("H56").Value = 31007.6
$ws.Range("I56").Value = 20346
$ws.Range("K56").Value = 20346
$ws.Range("M56").Value = -19655

# Row 61
$ws.Range("H61").Value = 2442.6365
$ws.Range("I61").Value = 2397.5789
$ws.Range("K61").Value = 2397.5789
$ws.Range("M61").Value = -2195.5789

# Row 100
$ws.Range("H100").Value = 4062.25
$ws.Range("I100").Value = 4566.5
$ws.Range("K100").Value = 4566.5
$ws.Range("M100").Value = -4025.5

# Row 113
$ws.Range("H113").Value = 2442.6365
$ws.Range("I113").Value = 2397.5789
$ws.Range("K113").Value = 2397.5789
$ws.Range("M113").Value = -227.5789

# Row 138
$ws.Range("H138").Value = 81770.836
$ws.Range("J138").Value = 81770.836
$ws.Range("L138").Value = 81770.836
$ws.Range("N138").Value = -92050.836

$ws = $wb.Sheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7721.8887
$ws.Range("I132").Value = 8356.714
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 25070.142
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -22540.142
$ws.Range("N132").Value = -21560

# Row 136
$ws.Range("H136").Value = 52635660
$ws.Range("I136").Value = 66667904
$ws.Range("K136").Value = 200003712
$ws.Range("M136").Value = -200001162
